$wb = $excel.ActiveWorkbook
$wsExample = $wb.Worksheets.Item("example")
$wsAll = $wb.Worksheets.Item("all")

# ---------------------------------------------------------------------------
# Add the four new "CostInvest" rows to the "all" sheet (rows 6-9)
# ---------------------------------------------------------------------------
$types = @("PowerPlants", "PowerPlants", "PowerPlants", "PowerPlants")
$techs = @("E_PV_DIST_RES", "E_SCO2", "E_OCAES", "E_BECCS")
$avgs  = @(1884, 2474, 1457, 6874)

for ($i = 0; $i -lt 4; $i++) {
    $r = 6 + $i
    $wsAll.Range("A$r").Value = "PowerPlants"
    $wsAll.Range("B$r").Value = "CostInvest"
    $wsAll.Range("C$r").Value = $techs[$i]
    $wsAll.Range("D$r").Value = "[M`$/GW]"
    $wsAll.Range("E$r").Value = $avgs[$i]
    $wsAll.Range("F$r").Value = "uniform"
}

# Row 6 gets its own (non-shared) formulas
$wsAll.Range("G6").Formula = "=E6*0.8"
$wsAll.Range("H6").Formula = "=E6*1.2"

# Rows 7-9 share one formula definition each
$wsAll.Range("G7:G9").Formula = "=E7*0.8"
$wsAll.Range("H7:H9").Formula = "=E7*1.2"

# Highlight the G6/H6 cells with a yellow fill
$wsAll.Range("G6:H6").Interior.Color = 65535

# E9 picks up the same number-format style as E5 (fontId 2 / style index 5)
$wsAll.Range("E5").Copy()
$wsAll.Range("E9").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Column width tweak on column F of "all"
# ---------------------------------------------------------------------------
$wsAll.Columns.Item(6).ColumnWidth = 16.75

# ---------------------------------------------------------------------------
# View/selection changes
# ---------------------------------------------------------------------------
$wsExample.Activate()
$excel.ActiveWindow.ScrollRow = 10
[void]$wsExample.Range("D8").Select()

$wsAll.Activate()
[void]$wsAll.Range("F5").Select()
